# Applies the "transp files and bfpat" edit to the RTMF workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) About sheet: update the "date updated" cell C1 (date serial)
# ---------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 44944

# ---------------------------------------------------------------
# 2) all_csv_SYVbT-passenger sheet: update raw vehicle-count data
#    in columns C (rail) and F (ships) for a number of state rows.
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("all_csv_SYVbT-passenger")

$wsData.Range("C17").Value = 3

$wsData.Range("C23").Value = 24

$wsData.Range("C29").Value = 450
$wsData.Range("F29").Value = 27

$wsData.Range("C35").Value = 72
$wsData.Range("F35").Value = 0

$wsData.Range("C41").Value = 0
$wsData.Range("F41").Value = 5

$wsData.Range("C47").Value = 0
$wsData.Range("F47").Value = 0

$wsData.Range("C53").Value = 29
$wsData.Range("F53").Value = 11

$wsData.Range("C59").Value = 41

$wsData.Range("C83").Value = 279
$wsData.Range("F83").Value = 21

$wsData.Range("C89").Value = 17

$wsData.Range("C107").Value = 31

$wsData.Range("C113").Value = 166
$wsData.Range("F113").Value = 43

$wsData.Range("C119").Value = 50
$wsData.Range("F119").Value = 5

# Row 125: both C125 and F125 become blank cells entirely
$wsData.Range("C125").ClearContents()
$wsData.Range("F125").ClearContents()

$wsData.Range("C131").Value = 5
$wsData.Range("F131").Value = 0

$wsData.Range("C137").Value = 26
$wsData.Range("F137").Value = 1

$wsData.Range("C143").Value = 32
$wsData.Range("F143").Value = 0

$wsData.Range("C161").Value = 20
$wsData.Range("F161").Value = 0

$wsData.Range("C185").Value = 190
$wsData.Range("F185").Value = 25

$wsData.Range("F191").Value = 7

$wsData.Range("C203").Value = 835
$wsData.Range("F203").Value = 17

$wsData.Range("C209").Value = 24

$wsData.Range("C215").Value = 4

$wsData.Range("C221").Value = 72
$wsData.Range("F221").Value = 1

$wsData.Range("C227").Value = 297
$wsData.Range("F227").Value = 0

$wsData.Range("C233").Value = 0
$wsData.Range("F233").Value = 0

$wsData.Range("C251").Value = 2
$wsData.Range("F251").Value = 4

$wsData.Range("C257").Value = 94
$wsData.Range("F257").Value = 19

$wsData.Range("C263").Value = 37
$wsData.Range("F263").Value = 6

$wsData.Range("C269").Value = 19
$wsData.Range("F269").Value = 19

$wsData.Range("C281").Value = 40
$wsData.Range("F281").Value = 2

$wsData.Range("C287").Value = 4
$wsData.Range("F287").Value = 0

# ---------------------------------------------------------------
# 3) Force a full recalculation so the dependent formulas on
#    "State RTMF" (B1, F5) and "RTMF-passengers" (E2, I2) pick up
#    the new raw data.
# ---------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------
# 4) Sheet view / selection bookkeeping to match the saved state.
# ---------------------------------------------------------------

# "State RTMF": move the cursor from G6 to F5
$wsState = $wb.Worksheets.Item("State RTMF")
$wsState.Range("F5").Select()

# "RTMF-freight": this sheet was tabSelected/active before; it no
# longer should be, so leave it without explicitly reselecting it.

# "all_csv_SYVbT-passenger" becomes the active/selected sheet with
# its entire grid selected.
$wsData.Activate()
$wsData.Range("A1:XFD1048576").Select()
